# Adding a new fluorophore column (mScarlet-I) to the properties table.
#
# The new data lives in a brand-new column W; every existing column from
# W onward shifts one place to the right (AH -> AI, dimension A1:AH35 ->
# A1:AI35, etc.). Inserting a whole column reproduces exactly that shift
# for every row in one shot, including formulas/styles, and keeps
# everything else (row spans, cached formula values, sharedStrings reuse)
# consistent automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh column before W -- shifts W:AH one column right to X:AI.
$ws.Columns("W").Insert()

# Populate the new column with mScarlet-I's data (header + 5 measurements).
$ws.Range("W7").Value = "mScarlet-I"
$ws.Range("W8").Value = 0.54
$ws.Range("W9").Value = 36
$ws.Range("W10").Value = 3.1
$ws.Range("W11").Value = 104000
$ws.Range("W12").Value = 225

# Leave the selection on the new data, matching the author's last position.
$ws.Range("W13").Select()
